$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws 'D2' '28.012.30'
Set-TextValue $ws 'E2' '  -2.07%  '
Set-TextValue $ws 'D3' '1.830.24'
Set-TextValue $ws 'E3' '  -1.08%  '
Set-TextValue $ws 'E4' '  -0.03%  '
Set-TextValue $ws 'D5' '325.58'
Set-TextValue $ws 'E5' '  -3.07%  '
Set-TextValue $ws 'E6' '  -0.11%  '
Set-TextValue $ws 'D7' '0.4638'
Set-TextValue $ws 'E7' '  -0.29%  '
Set-TextValue $ws 'D8' '0.3873'
Set-TextValue $ws 'E8' '  -1.13%  '
Set-TextValue $ws 'D9' '0.07868'
Set-TextValue $ws 'E9' '  -0.26%  '
Set-TextValue $ws 'D10' '0.9594'
Set-TextValue $ws 'E10' '  -2.43%  '
Set-TextValue $ws 'E11' '  -1.62%  '
Set-TextValue $ws 'D12' '1.881.53'
Set-TextValue $ws 'E12' '  -0.01%  '
Set-TextValue $ws 'D13' '5.669'
Set-TextValue $ws 'E13' '  -3.06%  '
Set-TextValue $ws 'E14' '  -1.60%  '
Set-TextValue $ws 'D15' '0.06769'
Set-TextValue $ws 'E15' '  -0.83%  '
Set-TextValue $ws 'D16' '87.17'
Set-TextValue $ws 'E16' '  -0.51%  '
Set-TextValue $ws 'D17' '1.001'
Set-TextValue $ws 'E17' '  -0.21%  '
Set-TextValue $ws 'D18' '0.000009921'
Set-TextValue $ws 'E18' '  -1.92%  '
Set-TextValue $ws 'D19' '16.60'
Set-TextValue $ws 'E19' '  -2.46%  '
Set-TextValue $ws 'E20' '  -0.11%  '
Set-TextValue $ws 'D21' '28.034.38'
Set-TextValue $ws 'E21' '  -2.02%  '
Set-TextValue $ws 'D22' '5.314'
Set-TextValue $ws 'E22' '  -1.56%  '
Set-TextValue $ws 'E23' '  -2.59%  '
Set-TextValue $ws 'D24' '2.094'
Set-TextValue $ws 'E24' '  -1.38%  '
Set-TextValue $ws 'D25' '2.049.62'
Set-TextValue $ws 'E25' '  -2.56%  '
Set-TextValue $ws 'E26' '  +0.31%  '
Set-TextValue $ws 'D27' '19.17'
Set-TextValue $ws 'E27' '  -1.19%  '
Set-TextValue $ws 'E28' '  -7.60%  '
Set-TextValue $ws 'E29' '  -2.32%  '
Set-TextValue $ws 'D30' '117.30'
Set-TextValue $ws 'E30' '  -0.16%  '
Set-TextValue $ws 'D31' '0.9362'
Set-TextValue $ws 'E31' '  -3.98%  '
Set-TextValue $ws 'D32' '0.09263'
Set-TextValue $ws 'E32' '  -1.87%  '
Set-TextValue $ws 'E33' '  -1.50%  '
Set-TextValue $ws 'D34' '1.318'
Set-TextValue $ws 'E34' '  -2.47%  '
Set-TextValue $ws 'D35' '3.286'
Set-TextValue $ws 'E35' '  -6.18%  '
Set-TextValue $ws 'E36' '  -4.33%  '
Set-TextValue $ws 'D37' '0.02141'
Set-TextValue $ws 'E37' '  -2.41%  '
Set-TextValue $ws 'E38' '  -1.21%  '
Set-TextValue $ws 'D39' '7.759'
Set-TextValue $ws 'E39' '  +2.36%  '
Set-TextValue $ws 'D40' '0.5580'
Set-TextValue $ws 'E40' '  -2.01%  '
Set-TextValue $ws 'D41' '9.874'
Set-TextValue $ws 'E41' '  -2.39%  '
Set-TextValue $ws 'D42' '0.1762'
Set-TextValue $ws 'E42' '  -1.51%  '
Set-TextValue $ws 'D43' '11.66'
Set-TextValue $ws 'E43' '  -1.44%  '
Set-TextValue $ws 'E44' '  -2.04%  '
Set-TextValue $ws 'D45' '0.07004'
Set-TextValue $ws 'E45' '  -1.96%  '
Set-TextValue $ws 'B46' 'RenderToken'
Set-TextValue $ws 'C46' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws 'D46' '2.122'
Set-TextValue $ws 'E46' '  -11.18%  '
Set-TextValue $ws 'B47' 'NEARProtocol'
Set-TextValue $ws 'C47' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws 'D47' '1.833'
Set-TextValue $ws 'E47' '  -3.90%  '
Set-TextValue $ws 'B48' 'Quant'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws 'D48' '112.91'
Set-TextValue $ws 'E48' '  -0.27%  '
Set-TextValue $ws 'B49' 'WEMIXToken'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws 'D49' '1.102'
Set-TextValue $ws 'E49' '  -12.15%  '
Set-TextValue $ws 'E50' '  -0.07%  '
Set-TextValue $ws 'E51' '  +0.47%  '
